$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.389.78'
$ws.Range("E2").Value = '  +2.23%  '
$ws.Range("D3").Value = '1.666.82'
$ws.Range("E3").Value = '  +1.38%  '
$ws.Range("E4").Value = '  -0.63%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '220.09'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.41%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.504'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.25%  '
$ws.Range("E7").Value = '  -0.60%  '
$ws.Range("E8").Value = '  +1.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0628'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.92'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +3.82%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0849'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.53%  '
$ws.Range("D12").Value = '1.899.59'
$ws.Range("E12").Value = '  +1.38%  '
$ws.Range("D13").Value = '1.668.00'
$ws.Range("E13").Value = '  +1.49%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.21'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.28%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.534'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.01%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '67.23'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +4.04%  '
$ws.Range("D17").Value = '27.372.54'
$ws.Range("E17").Value = '  +2.19%  '
$ws.Range("D18").Value = '0.0₃0738'
$ws.Range("E18").Value = '  +0.51%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '223.91'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +4.55%  '
$ws.Range("E20").Value = '  -0.72%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.82'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +9.19%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.45'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +1.28%  '
$ws.Range("E23").Value = '  +2.80%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.29'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.23%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.66'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.21%  '
$ws.Range("E26").Value = '  -0.68%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.44'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +3.67%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.119'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.86%  '
$ws.Range("E29").Value = '  +2.57%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0514'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +1.10%  '
$ws.Range("E31").Value = '  +1.48%  '
$ws.Range("E32").Value = '  +0.89%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.01'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.65%  '
$ws.Range("E34").Value = '  +2.54%  '
$ws.Range("D35").Value = '1.269.94'
$ws.Range("E35").Value = '  -1.02%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.45'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.09%  '
$ws.Range("E37").Value = '  -0.75%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.538'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.38%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.831'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.50%  '
$ws.Range("E40").Value = '  -0.56%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.812'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.88%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.39'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +1.70%  '
$ws.Range("D43").Value = '1.811.71'
$ws.Range("E43").Value = '  +1.58%  '
$ws.Range("E44").Value = '  -4.73%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '62.08'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +1.10%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '92.67'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.83%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.62'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +1.05%  '
$ws.Range("E48").Value = '  +0.10%  '
$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0984'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +1.39%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.66'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.17%  '
$ws.Range("E51").Value = '  +0.10%  '
